$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block (rows 6-19) holds one weekly pair of quality-grade rows
# (Primera/Segunda) per market date, oldest first. A new, most-recent week
# (2021-11-25, serial 44525) is being added to the series right after the
# two oldest fixed rows (2-5), so insert two new rows at row 6 and push the
# existing history down by two rows.
$ws.Rows("6:7").Insert()

# Fill in the new week's data (mirrors the existing Primera/Segunda rows).
$ws.Range("A6").Value2 = 11
$ws.Range("B6").Value2 = "Vega Monumental Concepción"
$ws.Range("C6").Value2 = "Bíobío"
$ws.Range("D6").Value2 = 44525
$ws.Range("E6").Value2 = 8
$ws.Range("F6").Value2 = 100112037
$ws.Range("G6").Value2 = "Cebollín"
$ws.Range("H6").Value2 = "Sin especificar"
$ws.Range("I6").Value2 = "Primera"
$ws.Range("J6").Value2 = 200
$ws.Range("K6").Value2 = 600
$ws.Range("L6").Value2 = 700
$ws.Range("M6").Value2 = 650
$ws.Range("N6").Value2 = "$/paquete 6 unidades"
$ws.Range("O6").Value2 = "Región de Ñuble"
$ws.Range("P6").Value2 = 108
$ws.Range("Q6").Value2 = 6
$ws.Range("R6").Value2 = "Hortaliza"

$ws.Range("A7").Value2 = 11
$ws.Range("B7").Value2 = "Vega Monumental Concepción"
$ws.Range("C7").Value2 = "Bíobío"
$ws.Range("D7").Value2 = 44525
$ws.Range("E7").Value2 = 8
$ws.Range("F7").Value2 = 100112037
$ws.Range("G7").Value2 = "Cebollín"
$ws.Range("H7").Value2 = "Sin especificar"
$ws.Range("I7").Value2 = "Segunda"
$ws.Range("J7").Value2 = 100
$ws.Range("K7").Value2 = 500
$ws.Range("L7").Value2 = 500
$ws.Range("M7").Value2 = 500
$ws.Range("N7").Value2 = "$/paquete 6 unidades"
$ws.Range("O7").Value2 = "Región de Ñuble"
$ws.Range("P7").Value2 = 83
$ws.Range("Q7").Value2 = 6
$ws.Range("R7").Value2 = "Hortaliza"
